$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) ANA's balance drops from 67728.23 to 43728.89
$ws.Cells.Item(3, 3).Value = 43728.89

# 2) New row for FRANCISCO right after ANA (row 4)
$ws.Rows(4).EntireRow.Insert()
$ws.Cells.Item(4, 1).NumberFormat = "@"
$ws.Cells.Item(4, 1).Value = "004567324"
$ws.Cells.Item(4, 1).ClearFormats()
$ws.Cells.Item(4, 2).Value = "FRANCISCO"
$ws.Cells.Item(4, 3).Value = 34627.71

# 3) Three new rows (DANIEL, RICARDO, LUIZ) inserted right before the HFR
#    row (which, after the insert above, now sits at row 8)
$ws.Rows(8).EntireRow.Insert()
$ws.Rows(8).EntireRow.Insert()
$ws.Rows(8).EntireRow.Insert()

$ws.Cells.Item(8, 1).NumberFormat = "@"
$ws.Cells.Item(8, 1).Value = "004493324"
$ws.Cells.Item(8, 1).ClearFormats()
$ws.Cells.Item(8, 2).Value = "DANIEL"
$ws.Cells.Item(8, 3).Value = 7695.34

$ws.Cells.Item(9, 1).NumberFormat = "@"
$ws.Cells.Item(9, 1).Value = "004505474"
$ws.Cells.Item(9, 1).ClearFormats()
$ws.Cells.Item(9, 2).Value = "RICARDO"
$ws.Cells.Item(9, 3).Value = 7551.89

$ws.Cells.Item(10, 1).NumberFormat = "@"
$ws.Cells.Item(10, 1).Value = "004298042"
$ws.Cells.Item(10, 1).ClearFormats()
$ws.Cells.Item(10, 2).Value = "LUIZ"
$ws.Cells.Item(10, 3).Value = 7169.58

# 4) New row for PAULA right before the CARLOS row (now at row 12)
$ws.Rows(12).EntireRow.Insert()
$ws.Cells.Item(12, 1).NumberFormat = "@"
$ws.Cells.Item(12, 1).Value = "004471893"
$ws.Cells.Item(12, 1).ClearFormats()
$ws.Cells.Item(12, 2).Value = "PAULA"
$ws.Cells.Item(12, 3).Value = 5121.13

# 5) Remove the old DANIEL row (balance 143.45) that used to sit between
#    ALEXANDRE (147.18) and KARINA (137.66). After the five rows inserted
#    above, it now lives at row 90.
$ws.Rows(90).EntireRow.Delete()
